# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages update).

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows keyed by their F-column (want-to-go count) cell
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6767
$wsExhibit.Range("F13").Value = 412
$wsExhibit.Range("F17").Value = 3397
$wsExhibit.Range("F21").Value = 2032
$wsExhibit.Range("F28").Value = 136

# Sheet "全部类型" - same events, shifted row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6767
$wsAll.Range("F14").Value = 412
$wsAll.Range("F18").Value = 3397
$wsAll.Range("F22").Value = 2032
$wsAll.Range("F29").Value = 136
